$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.530.66"
$ws.Range("E2").Value = "  +0.61%  "

# Row 3
$ws.Range("D3").Value = "1.872.58"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "

# Row 6
$ws.Range("E6").Value = "  +0.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4732"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2908"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.74%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06470"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.10%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.12%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07699"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.04%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7407"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.44%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.70%  "

# Row 14
$ws.Range("D14").Value = "1.870.59"
$ws.Range("E14").Value = "  -0.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.163"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.63%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.79%  "

# Row 17
$ws.Range("D17").Value = "30.596.78"
$ws.Range("E17").Value = "  +0.89%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.63%  "

# Row 19
$ws.Range("E19").Value = "  +0.04%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007482"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.78%  "

# Row 21
$ws.Range("D21").Value = "2.119.28"
$ws.Range("E21").Value = "  +0.30%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.255"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.27%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.174"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.19%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.197"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.76%  "

# Row 26
$ws.Range("E26").Value = "  -0.19%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.907"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.69%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09984"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.93%  "

# Row 30
$ws.Range("E30").Value = "  -2.63%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.513"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.255"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.49%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.085"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.79%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04792"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.36%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.119"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.32%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6937"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.54%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01855"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.47%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.756"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.21%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.238"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.96%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.69%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.967"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.81%  "

# Row 43
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.15%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4164"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.66%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8338"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.24%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.81%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.368"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "

# Row 48
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.31%  "

# Row 49
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.987"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.45%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "913.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.49%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05657"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.60%  "
